# Fill in the "floor type" lookup table (Table2, G1:H4) with the three
# flooring options and their cost per square foot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Hardwood"
$ws.Range("H2").Value = 1.99
$ws.Range("G3").Value = "Carpet"
$ws.Range("H3").Value = 3.99
$ws.Range("G4").Value = "Tile"
$ws.Range("H4").Value = 4.99

# Row 3: Hardwood, 12 x 20
$ws.Range("A3").Value = "Hardwood"
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 20
$ws.Range("D3").Formula = '=($B$3 * $C$3) * H2'

# Row 4: Carpet, 12 x 20
$ws.Range("A4").Value = "Carpet"
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 20
$ws.Range("D4").Formula = '=($B$4 * $C$4) * H3'

# Row 5: Tile, 12 x 20
$ws.Range("A5").Value = "Tile"
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 20
$ws.Range("D5").Formula = '=($B$5* $C$5) * H4'

# Row 6: Hardwood, 5 x 9
$ws.Range("A6").Value = "Hardwood"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 9
$ws.Range("D6").Formula = '=($B$6 * $C$6) * H2'

# Row 7: Carpet, 5 x 9
$ws.Range("A7").Value = "Carpet"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 9
$ws.Range("D7").Formula = '=($B$7 * $C$7) * H3'

# Row 8: overall total cost
$ws.Range("D8").Formula = '=($D$3 + $D$4 + $D$5 + $D$6 + $D$7)'

# Leave the selection where the author left it
$ws.Range("E12").Select()
